$d = $word.ActiveDocument

# Locate the full old sentence within the document.
$old = "Stock Market Trend Prediction: Apply Gradient Boosting Machines (GBM) to build a model for predicting stock market trends and making investment decisions."
$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target sentence to replace."
}

$startPos = $rng.Start

# New wording, split into the pieces the edit introduces; each piece ends up
# as its own run (matching the target markup) even though all runs share
# identical (strikethrough) formatting.
$parts = @(
    "Classifying and Predicting Stock Market States",
    ": Apply Gradient Boosting Machines (GBM) to build a model ",
    "with Hidden Markov Model (HMM) ",
    "for predicting stock market trends and making investment decisions."
)
$n = $parts.Length

# Put the LAST piece into the existing (anchor) paragraph, reusing its
# original run/paragraph so the paragraph keeps its identity (paraId, rsids,
# etc. - none of which the diff touches).
$rng.Text = $parts[$n - 1]
$anchor = $d.Range($startPos, $startPos).Paragraphs(1)
$anchor.Range.Font.StrikeThrough = 1

# For each remaining piece (from second-to-last back to the first), insert a
# new empty paragraph immediately before the anchor paragraph, fill it with
# that piece of text as its own run, then merge it back into the anchor by
# deleting the pilcrow that separates them. Deleting a pilcrow keeps the
# properties of the paragraph AFTER it (the anchor), so the anchor's
# original identity survives every merge, while every inserted piece still
# keeps being written out as a separate <w:r>.
for ($k = $n - 2; $k -ge 0; $k--) {
    $insParaPoint = $d.Range($startPos, $startPos)
    $insParaPoint.InsertParagraphBefore()

    $segRange = $d.Range($startPos, $startPos)
    $segRange.InsertAfter($parts[$k])
    $segRange.Font.StrikeThrough = 1

    $pilcrowPos = $segRange.End
    $d.Range($pilcrowPos, $pilcrowPos + 1).Delete()
}
